$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prepare formats for the 5 new rows (85-89) by copying the format of the
# last existing data row (row 84), then fill in the values.
$ws.Range("A84:V84").Copy()
$ws.Range("A85:V89").PasteSpecial(-4122)

# Row 85 -> Indice 84 : Igman K. 0 x 1 Borac Banja Luka
$ws.Cells.Item(85, 1).Value = 84
$ws.Cells.Item(85, 2).Value = "bosnia-and-herzegovina"
$ws.Cells.Item(85, 3).Value = "premijer-liga-bih"
$ws.Cells.Item(85, 4).Value = "2023-2024"
$ws.Cells.Item(85, 5).Value = 45255.54166666666
$ws.Cells.Item(85, 6).Value = "Igman K."
$ws.Cells.Item(85, 7).Value = 0
$ws.Cells.Item(85, 8).Value = "Borac Banja Luka"
$ws.Cells.Item(85, 9).Value = 1
$ws.Cells.Item(85, 10).Value = 4.77
$ws.Cells.Item(85, 11).Value = "24/11/2023 01:12"
$ws.Cells.Item(85, 12).Value = 4.38
$ws.Cells.Item(85, 13).Value = "25/11/2023 12:58"
$ws.Cells.Item(85, 14).Value = 3.52
$ws.Cells.Item(85, 15).Value = "24/11/2023 01:12"
$ws.Cells.Item(85, 16).Value = 4.15
$ws.Cells.Item(85, 17).Value = "25/11/2023 12:58"
$ws.Cells.Item(85, 18).Value = 1.63
$ws.Cells.Item(85, 19).Value = "24/11/2023 01:12"
$ws.Cells.Item(85, 20).Value = 1.66
$ws.Cells.Item(85, 21).Value = "25/11/2023 12:58"
$ws.Cells.Item(85, 22).Value = "https://www.betexplorer.com/football/bosnia-and-herzegovina/premijer-liga-bih/igman-konjic-borac-banja-luka/KlXFGV7N/"

# Row 86 -> Indice 85 : Sloga Doboj 1 x 0 Zvijezda 09
$ws.Cells.Item(86, 1).Value = 85
$ws.Cells.Item(86, 2).Value = "bosnia-and-herzegovina"
$ws.Cells.Item(86, 3).Value = "premijer-liga-bih"
$ws.Cells.Item(86, 4).Value = "2023-2024"
$ws.Cells.Item(86, 5).Value = 45256.54166666666
$ws.Cells.Item(86, 6).Value = "Sloga Doboj"
$ws.Cells.Item(86, 7).Value = 1
$ws.Cells.Item(86, 8).Value = "Zvijezda 09"
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 10).Value = 1.47
$ws.Cells.Item(86, 11).Value = "25/11/2023 07:42"
$ws.Cells.Item(86, 12).Value = 1.58
$ws.Cells.Item(86, 13).Value = "26/11/2023 12:57"
$ws.Cells.Item(86, 14).Value = 4.06
$ws.Cells.Item(86, 15).Value = "25/11/2023 07:42"
$ws.Cells.Item(86, 16).Value = 3.57
$ws.Cells.Item(86, 17).Value = "26/11/2023 12:58"
$ws.Cells.Item(86, 18).Value = 5.28
$ws.Cells.Item(86, 19).Value = "25/11/2023 07:42"
$ws.Cells.Item(86, 20).Value = 6.32
$ws.Cells.Item(86, 21).Value = "26/11/2023 12:57"
$ws.Cells.Item(86, 22).Value = "https://www.betexplorer.com/football/bosnia-and-herzegovina/premijer-liga-bih/sloga-doboj-zvijezda-09/vZDnB9Vp/"

# Row 87 -> Indice 86 : Zeljeznicar 1 x 1 Posusje
$ws.Cells.Item(87, 1).Value = 86
$ws.Cells.Item(87, 2).Value = "bosnia-and-herzegovina"
$ws.Cells.Item(87, 3).Value = "premijer-liga-bih"
$ws.Cells.Item(87, 4).Value = "2023-2024"
$ws.Cells.Item(87, 5).Value = 45256.54166666666
$ws.Cells.Item(87, 6).Value = "Zeljeznicar"
$ws.Cells.Item(87, 7).Value = 1
$ws.Cells.Item(87, 8).Value = "Posusje"
$ws.Cells.Item(87, 9).Value = 1
$ws.Cells.Item(87, 10).Value = 1.74
$ws.Cells.Item(87, 11).Value = "25/11/2023 07:42"
$ws.Cells.Item(87, 12).Value = 1.81
$ws.Cells.Item(87, 13).Value = "26/11/2023 12:54"
$ws.Cells.Item(87, 14).Value = 3.29
$ws.Cells.Item(87, 15).Value = "25/11/2023 07:42"
$ws.Cells.Item(87, 16).Value = 3.36
$ws.Cells.Item(87, 17).Value = "26/11/2023 12:57"
$ws.Cells.Item(87, 18).Value = 4.25
$ws.Cells.Item(87, 19).Value = "25/11/2023 07:42"
$ws.Cells.Item(87, 20).Value = 4.5
$ws.Cells.Item(87, 21).Value = "26/11/2023 12:57"
$ws.Cells.Item(87, 22).Value = "https://www.betexplorer.com/football/bosnia-and-herzegovina/premijer-liga-bih/zeljeznicar-posusje/2DMKFkNT/"

# Row 88 -> Indice 87 : GOSK Gabela 1 x 5 Zrinjski
$ws.Cells.Item(88, 1).Value = 87
$ws.Cells.Item(88, 2).Value = "bosnia-and-herzegovina"
$ws.Cells.Item(88, 3).Value = "premijer-liga-bih"
$ws.Cells.Item(88, 4).Value = "2023-2024"
$ws.Cells.Item(88, 5).Value = 45256.625
$ws.Cells.Item(88, 6).Value = "GOSK Gabela"
$ws.Cells.Item(88, 7).Value = 1
$ws.Cells.Item(88, 8).Value = "Zrinjski"
$ws.Cells.Item(88, 9).Value = 5
$ws.Cells.Item(88, 10).Value = 5.22
$ws.Cells.Item(88, 11).Value = "25/11/2023 03:12"
$ws.Cells.Item(88, 12).Value = 9.68
$ws.Cells.Item(88, 13).Value = "26/11/2023 14:57"
$ws.Cells.Item(88, 14).Value = 3.89
$ws.Cells.Item(88, 15).Value = "25/11/2023 03:12"
$ws.Cells.Item(88, 16).Value = 5.03
$ws.Cells.Item(88, 17).Value = "26/11/2023 14:57"
$ws.Cells.Item(88, 18).Value = 1.5
$ws.Cells.Item(88, 19).Value = "25/11/2023 03:12"
$ws.Cells.Item(88, 20).Value = 1.3
$ws.Cells.Item(88, 21).Value = "26/11/2023 14:57"
$ws.Cells.Item(88, 22).Value = "https://www.betexplorer.com/football/bosnia-and-herzegovina/premijer-liga-bih/nk-gosk-gabela-zrinjski/WS0YYCNi/"

# Row 89 -> Indice 88 : Siroki Brijeg 2 x 2 Velez Mostar
$ws.Cells.Item(89, 1).Value = 88
$ws.Cells.Item(89, 2).Value = "bosnia-and-herzegovina"
$ws.Cells.Item(89, 3).Value = "premijer-liga-bih"
$ws.Cells.Item(89, 4).Value = "2023-2024"
$ws.Cells.Item(89, 5).Value = 45256.82291666666
$ws.Cells.Item(89, 6).Value = "Siroki Brijeg"
$ws.Cells.Item(89, 7).Value = 2
$ws.Cells.Item(89, 8).Value = "Velez Mostar"
$ws.Cells.Item(89, 9).Value = 2
$ws.Cells.Item(89, 10).Value = 2.97
$ws.Cells.Item(89, 11).Value = "25/11/2023 08:12"
$ws.Cells.Item(89, 12).Value = 3.05
$ws.Cells.Item(89, 13).Value = "26/11/2023 19:44"
$ws.Cells.Item(89, 14).Value = 3.04
$ws.Cells.Item(89, 15).Value = "25/11/2023 08:12"
$ws.Cells.Item(89, 16).Value = 3.09
$ws.Cells.Item(89, 17).Value = "26/11/2023 19:44"
$ws.Cells.Item(89, 18).Value = 2.27
$ws.Cells.Item(89, 19).Value = "25/11/2023 08:12"
$ws.Cells.Item(89, 20).Value = 2.38
$ws.Cells.Item(89, 21).Value = "26/11/2023 19:44"
$ws.Cells.Item(89, 22).Value = "https://www.betexplorer.com/football/bosnia-and-herzegovina/premijer-liga-bih/siroki-brijeg-velez-mostar/hx2UZh8o/"
